# Fix grammar and spelling errors in the "Getting the team together" bullet:
#   ". Our team consists four people."                              -> "... consists of four people."
#   "... a specific role in the realization the project;"           -> "... realization of the project;"
#
# Both fixes insert the missing word "of" at a precise point in the middle of
# an existing run. We do this with a zero-length Range + Text assignment
# (the COM equivalent of clicking at that spot and typing " of"), which is
# how the original author's edit naturally produced three runs out of one
# (the untouched head, the newly typed " of", and the untouched tail).

$d = $word.ActiveDocument

function Insert-Of($afterPos) {
    # Type " of" at the collapsed insertion point.
    $ins = $d.Range($afterPos, $afterPos)
    $ins.Text = " of"

    # Nudge a character property on the text we just typed and immediately
    # revert it. This keeps the newly typed text as its own run(s) instead
    # of the engine quietly re-coalescing it back into the neighbouring
    # runs just because the resolved formatting happens to be identical.
    $touched = $d.Range($afterPos, $afterPos + 3)
    $touched.Font.Bold = $true
    $touched.Font.Bold = $false
}

# --- Edit 1: "Our team consists four people" -> "Our team consists of four people" ---
$rng1 = $d.Content
$rng1.Find.Execute("consists", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Insert-Of $rng1.End

# --- Edit 2: "realization the project" -> "realization of the project" ---
# "realization" occurs several times in the document, so first narrow down
# to the unique sentence before locating the exact insertion point inside it.
$scope = $d.Content
$scope.Find.Execute("role in the realization the project", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2 = $d.Range($scope.Start, $scope.End)
$rng2.Find.Execute("realization", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Insert-Of $rng2.End
